$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row data: Row number, new B (Coin), new C (Link), new D (Price), new E (Volume 1h)
# Only rows whose values actually changed are listed; blank (empty string) means "no change" for that column.
$updates = @(
    @{ Row = 2; B = ''; C = ''; D = '37.214.79'; E = '  +0.28%  ' },
    @{ Row = 3; B = ''; C = ''; D = '2.079.20'; E = '  -0.16%  ' },
    @{ Row = 4; B = ''; C = ''; D = ''; E = '  +0.19%  ' },
    @{ Row = 5; B = ''; C = ''; D = '252.78'; E = '  +1.64%  ' },
    @{ Row = 6; B = ''; C = ''; D = ''; E = '  +4.40%  ' },
    @{ Row = 7; B = ''; C = ''; D = '62.09'; E = '  +21.88%  ' },
    @{ Row = 8; B = ''; C = ''; D = '1.00'; E = '  +0.03%  ' },
    @{ Row = 9; B = ''; C = ''; D = '61.73'; E = '  +2.51%  ' },
    @{ Row = 10; B = ''; C = ''; D = ''; E = '  +5.87%  ' },
    @{ Row = 11; B = ''; C = ''; D = '0.0807'; E = '  +9.82%  ' },
    @{ Row = 12; B = ''; C = ''; D = ''; E = '  +2.73%  ' },
    @{ Row = 13; B = ''; C = ''; D = '15.66'; E = '  +3.11%  ' },
    @{ Row = 14; B = ''; C = ''; D = '2.386.77'; E = '  +0.14%  ' },
    @{ Row = 15; B = ''; C = ''; D = '0.829'; E = '  +0.33%  ' },
    @{ Row = 16; B = ''; C = ''; D = '5.47'; E = '  +8.33%  ' },
    @{ Row = 17; B = ''; C = ''; D = '2.085.42'; E = '  +0.11%  ' },
    @{ Row = 18; B = ''; C = ''; D = '37.230.70'; E = '  +0.69%  ' },
    @{ Row = 19; B = ''; C = ''; D = '74.80'; E = '  +3.99%  ' },
    @{ Row = 20; B = ''; C = ''; D = '0.0₃0926'; E = '  +13.44%  ' },
    @{ Row = 21; B = ''; C = ''; D = '15.04'; E = '  +14.44%  ' },
    @{ Row = 22; B = ''; C = ''; D = '5.47'; E = '  +5.56%  ' },
    @{ Row = 23; B = ''; C = ''; D = '240.48'; E = '  +1.15%  ' },
    @{ Row = 24; B = ''; C = ''; D = '0.999'; E = '  -0.15%  ' },
    @{ Row = 25; B = ''; C = ''; D = ''; E = '  -1.38%  ' },
    @{ Row = 26; B = ''; C = ''; D = '171.81'; E = '  +1.60%  ' },
    @{ Row = 27; B = 'PancakeSwap'; C = 'https://coinranking.com/coin/ncYFcP709+pancakeswap-cake'; D = '2.16'; E = '  +9.36%  ' },
    @{ Row = 28; B = 'Cosmos'; C = 'https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom'; D = '9.29'; E = '  +1.90%  ' },
    @{ Row = 29; B = 'EthereumClassic'; C = 'https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc'; D = '20.42'; E = '  -1.10%  ' },
    @{ Row = 30; B = ''; C = ''; D = ''; E = '  +4.00%  ' },
    @{ Row = 31; B = ''; C = ''; D = '4.84'; E = '  +8.92%  ' },
    @{ Row = 32; B = ''; C = ''; D = ''; E = '  +3.62%  ' },
    @{ Row = 33; B = ''; C = ''; D = '0.0638'; E = '  +5.64%  ' },
    @{ Row = 34; B = ''; C = ''; D = '4.46'; E = '  +10.86%  ' },
    @{ Row = 35; B = ''; C = ''; D = '0.0901'; E = '  +0.46%  ' },
    @{ Row = 36; B = 'BinanceUSD'; C = 'https://coinranking.com/coin/vSo2fu9iE1s0Y+binanceusd-busd'; D = '1.00'; E = '  +0.14%  ' },
    @{ Row = 37; B = 'LidoDAOToken'; C = 'https://coinranking.com/coin/Pe93bIOD2+lidodaotoken-ldo'; D = '2.32'; E = '  +2.44%  ' },
    @{ Row = 38; B = ''; C = ''; D = ''; E = '  -3.43%  ' },
    @{ Row = 39; B = ''; C = ''; D = ''; E = '  +25.13%  ' },
    @{ Row = 40; B = ''; C = ''; D = '1.36'; E = '  +3.92%  ' },
    @{ Row = 41; B = ''; C = ''; D = '18.59'; E = '  +6.33%  ' },
    @{ Row = 42; B = ''; C = ''; D = ''; E = '  +2.85%  ' },
    @{ Row = 43; B = ''; C = ''; D = '1.16'; E = '  +1.93%  ' },
    @{ Row = 44; B = ''; C = ''; D = '99.33'; E = '  +2.06%  ' },
    @{ Row = 45; B = ''; C = ''; D = '4.33'; E = '  +22.80%  ' },
    @{ Row = 46; B = ''; C = ''; D = ''; E = '  +1.44%  ' },
    @{ Row = 47; B = 'THORChain'; C = 'https://coinranking.com/coin/ybmU-kKU+thorchain-rune'; D = '4.65'; E = '  +18.15%  ' },
    @{ Row = 48; B = 'RenderToken'; C = 'https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr'; D = '2.54'; E = '  +12.73%  ' },
    @{ Row = 49; B = ''; C = ''; D = '1.312.56'; E = '  +0.83%  ' },
    @{ Row = 50; B = ''; C = ''; D = '2.95'; E = '  -0.68%  ' },
    @{ Row = 51; B = ''; C = ''; D = '6.95'; E = '  +1.45%  ' }
)


foreach ($item in $updates) {
    $r = $item.Row

    if ($item.B -ne '') {
        $ws.Cells.Item($r, 2).Value2 = $item.B
    }
    if ($item.C -ne '') {
        $ws.Cells.Item($r, 3).Value2 = $item.C
    }
    if ($item.D -ne '') {
        $cellD = $ws.Cells.Item($r, 4)
        $origStyle = $cellD.Style
        $cellD.NumberFormat = "@"
        $cellD.Value2 = $item.D
        $cellD.Style = $origStyle
    }
    if ($item.E -ne '') {
        $ws.Cells.Item($r, 5).Value2 = $item.E
    }
}
